# Clear the cells on "Register_invalid" that previously held the literal
# text "empty", turning them into genuinely blank cells. This removes the
# now-unused "empty" entry from the shared strings table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register_invalid")

$emptyCells = @("A2", "B2", "C2", "B3", "C3", "C4", "A5", "B6")
foreach ($cellRef in $emptyCells) {
    $ws.Range($cellRef).ClearContents()
}

# Activate the Register_invalid sheet and change the selected cell from
# E6 to B6, matching the new saved view state.
$ws.Activate()
$ws.Range("B6").Select()
